$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 152, pushing existing rows 152:243 down to 153:244
$ws.Rows("152:152").Insert()

# Populate the newly inserted row with the new data record
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C152").Value = "Ñuble"
$ws.Range("D152").Value = 44813
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100108
$ws.Range("H152").Value = "Tropicales y subtropicales"
$ws.Range("I152").Value = 100108005
$ws.Range("J152").Value = "Piña"
$ws.Range("K152").Value = "Caramelo"
$ws.Range("L152").Value = "Segunda"
$ws.Range("M152").Value = 120
$ws.Range("N152").Value = 20000
$ws.Range("O152").Value = 21000
$ws.Range("P152").Value = 20500
$ws.Range("Q152").Value = "$/caja 14 unidades"
$ws.Range("R152").Value = "Ecuador"
$ws.Range("S152").Value = 1464
$ws.Range("T152").Value = 14
